$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.573.52"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "3.527.05"
$ws.Range("E3").Value = "  -2.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.62%  "

$ws.Range("D7").Value = "3.524.62"
$ws.Range("E7").Value = "  -2.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.10%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.20%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.409"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.75%  "

$ws.Range("D13").Value = "4.127.95"
$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("E14").Value = "  -7.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.78%  "

$ws.Range("D16").Value = "3.520.43"
$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").Value = "66.443.35"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.590"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("D25").Value = "3.669.54"
$ws.Range("E25").Value = "  -2.23%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -7.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.21%  "

$ws.Range("E29").Value = "  -2.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").Value = "3.535.92"
$ws.Range("E32").Value = "  -1.86%  "

$ws.Range("E33").Value = "  -1.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.13%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  -9.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.26%  "

$ws.Range("E38").Value = "  -4.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "174.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0819"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.53%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.28%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.857"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.78%  "

$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("E49").Value = "  -6.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.906"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.36%  "
